# Auto-generated edit script applying numeric corrections to Kraken Profits leve-crafting sheets
# per scheduled runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 221.9
$ws.Range("I33").Value = 235.44444
$ws.Range("K33").Value = 235.44444
$ws.Range("M33").Value = -6.444439999999986
$ws.Range("H86").Value = 2624.3333
$ws.Range("I86").Value = 832.3333
$ws.Range("K86").Value = 832.3333
$ws.Range("M86").Value = 290.6667
$ws.Range("H89").Value = 2624.3333
$ws.Range("I89").Value = 832.3333
$ws.Range("K89").Value = 4161.6665
$ws.Range("M89").Value = 1454.3335
$ws.Range("H99").Value = 313.8
$ws.Range("I99").Value = 342.25
$ws.Range("J99").Value = 200
$ws.Range("K99").Value = 1026.75
$ws.Range("L99").Value = 600
$ws.Range("M99").Value = 471.25
$ws.Range("N99").Value = -3596
$ws.Range("H101").Value = 972.5
$ws.Range("J101").Value = 972.5
$ws.Range("L101").Value = 2917.5
$ws.Range("N101").Value = -6161.5
$ws.Range("H106").Value = 3996.5
$ws.Range("I106").Value = 3996.5
$ws.Range("K106").Value = 3996.5
$ws.Range("M106").Value = -3365.5
$ws.Range("H107").Value = 2313.8572
$ws.Range("I107").Value = 1539.6
$ws.Range("K107").Value = 1539.6
$ws.Range("M107").Value = 380.4000000000001
$ws.Range("H111").Value = 600
$ws.Range("I111").Value = 600
$ws.Range("K111").Value = 1800
$ws.Range("M111").Value = 1267
$ws.Range("H135").Value = 1178.3846
$ws.Range("I135").Value = 1178.3846
$ws.Range("K135").Value = 10605.4614
$ws.Range("M135").Value = -8070.4614
$ws.Range("H138").Value = 3698.5264
$ws.Range("I138").Value = 3047.3333
$ws.Range("K138").Value = 9141.999899999999
$ws.Range("M138").Value = -4001.999899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7329.6665
$ws.Range("I32").Value = 7058.5
$ws.Range("K32").Value = 7058.5
$ws.Range("M32").Value = -6771.5
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").ClearContents()
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = 0
$ws.Range("H61").Value = 3178.8
$ws.Range("I61").Value = 2977.3333
$ws.Range("K61").Value = 2977.3333
$ws.Range("M61").Value = -2765.3333
$ws.Range("H94").Value = 9999
$ws.Range("J94").Value = 9999
$ws.Range("L94").Value = 9999
$ws.Range("N94").Value = -11801
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").ClearContents()
$ws.Range("N109").Value = 0
$ws.Range("H136").Value = 3178.8
$ws.Range("I136").Value = 2977.3333
$ws.Range("K136").Value = 8931.999899999999
$ws.Range("M136").Value = -6381.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 1064
$ws.Range("J23").Value = 1064
$ws.Range("L23").Value = 1064
$ws.Range("N23").Value = -1630
$ws.Range("H29").Value = 8403.666999999999
$ws.Range("I29").Value = 105.5
$ws.Range("J29").Value = 25000
$ws.Range("K29").Value = 105.5
$ws.Range("L29").Value = 25000
$ws.Range("M29").Value = 183.5
$ws.Range("N29").Value = -25578
$ws.Range("H94").Value = 994.2222
$ws.Range("I94").Value = 991.3333
$ws.Range("K94").Value = 991.3333
$ws.Range("M94").Value = -540.3333
$ws.Range("H95").Value = 12666.667
$ws.Range("J95").Value = 12666.667
$ws.Range("L95").Value = 12666.667
$ws.Range("N95").Value = -18158.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H69").Value = 7750
$ws.Range("I69").Value = 7750
$ws.Range("K69").Value = 7750
$ws.Range("M69").Value = -7001
$ws.Range("H72").Value = 7750
$ws.Range("I72").Value = 7750
$ws.Range("K72").Value = 23250
$ws.Range("M72").Value = -19506
$ws.Range("H93").Value = 60000
$ws.Range("I93").Value = 60000
$ws.Range("K93").Value = 60000
$ws.Range("M93").Value = -58128
$ws.Range("H107").Value = 787.25
$ws.Range("I107").Value = 716.6667
$ws.Range("K107").Value = 716.6667
$ws.Range("M107").Value = 1203.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 20.75
$ws.Range("I33").Value = 22.666666
$ws.Range("J33").Value = 19.6
$ws.Range("K33").Value = 135.999996
$ws.Range("L33").Value = 117.6
$ws.Range("M33").Value = 147.000004
$ws.Range("N33").Value = -683.6
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("H97").Value = 795.3333
$ws.Range("J97").Value = 700
$ws.Range("L97").Value = 2100
$ws.Range("N97").Value = -3092
$ws.Range("H129").Value = 3000
$ws.Range("J129").Value = 3000
$ws.Range("L129").Value = 9000
$ws.Range("N129").Value = -19000

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 19722222
$ws.Range("I11").Value = 19722222
$ws.Range("K11").Value = 19722222
$ws.Range("M11").Value = -19722083
$ws.Range("H92").Value = 12449.75
$ws.Range("J92").Value = 12449.75
$ws.Range("L92").Value = 12449.75
$ws.Range("N92").Value = -16193.75
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").ClearContents()
$ws.Range("N98").Value = 0
$ws.Range("H113").Value = 5959.8
$ws.Range("I113").Value = 4900
$ws.Range("K113").Value = 4900
$ws.Range("M113").Value = -2730
$ws.Range("H132").Value = 2603.4285
$ws.Range("I132").Value = 2339.3333
$ws.Range("J132").Value = 4188
$ws.Range("K132").Value = 7017.999899999999
$ws.Range("L132").Value = 12564
$ws.Range("M132").Value = -4487.999899999999
$ws.Range("N132").Value = -17624

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4287
$ws.Range("I46").Value = 2500
$ws.Range("J46").Value = 4584.8335
$ws.Range("K46").Value = 2500
$ws.Range("L46").Value = 4584.8335
$ws.Range("M46").Value = -2312
$ws.Range("N46").Value = -4960.8335
$ws.Range("H111").Value = 30000
$ws.Range("J111").Value = 30000
$ws.Range("L111").Value = 30000
$ws.Range("N111").Value = -38180

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 20010
$ws.Range("I20").Value = 20010
$ws.Range("K20").Value = 20010
$ws.Range("M20").Value = -19770
$ws.Range("H100").Value = 2758.3076
$ws.Range("I100").Value = 2665.889
$ws.Range("K100").Value = 5331.778
$ws.Range("M100").Value = -4790.778
$ws.Range("H103").Value = 20000
$ws.Range("J103").Value = 20000
$ws.Range("L103").Value = 20000
$ws.Range("N103").Value = -22344
$ws.Range("H122").Value = 5550.3
$ws.Range("I122").Value = 5187.25
$ws.Range("J122").Value = 7002.5
$ws.Range("K122").Value = 15561.75
$ws.Range("L122").Value = 21007.5
$ws.Range("M122").Value = -13111.75
$ws.Range("N122").Value = -25907.5

